# Avances en notebook y documento
# Applies the edits described by the diff:
#  - Hoja1: replace GRUPO category rows (Grande/Mediana/Pequeña/Micro ->
#    Manufactura/Servicos/Comercio), update statistic/p-value numbers,
#    update the decision text for the now-significant Comercio row, and
#    remove the now-empty "Micro" row (shifting the Levene's-test block up).
#  - Hoja2 / Hoja3: tidy the P-VALOR cell formatting so it re-uses the
#    existing bordered+centered style instead of a duplicate one.
#  - Make Hoja1 the active/selected sheet (it was Hoja3 before).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Hoja1: update the normality-test table
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Hoja1")

$ws1.Range("A2").Value = "Manufactura"
$ws1.Range("B2").Value = 0.9226
$ws1.Range("C2").Value = 0.2399

$ws1.Range("A3").Value = "Servicos"
$ws1.Range("B3").Value = 0.9524
$ws1.Range("C3").Value = 0.7355

$ws1.Range("A4").Value = "Comercio"
$ws1.Range("B4").Value = 0.6173
$ws1.Range("C4").Value = 0.0007
$ws1.Range("D4").Value = "Se rechaza la hipótesis nula. Los datos no siguen una distribución normal."

# Remove the now-obsolete "Micro" row (row 5) and shift everything below
# it up by one, matching the new A1:D7 used range.
$ws1.Rows.Item(5).Delete() | Out-Null

# ---------------------------------------------------------------
# Hoja2 / Hoja3: re-apply the border+center formatting on the P-VALOR
# cell by copying it from an already-consistent cell, which lets the
# duplicate style definition collapse into the shared one.
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Hoja2")
$ws3 = $wb.Worksheets.Item("Hoja3")

$ws1.Range("B2").Copy() | Out-Null
$ws2.Range("B2").PasteSpecial(-4122) | Out-Null

$ws1.Range("B2").Copy() | Out-Null
$ws3.Range("B2").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------
# Make Hoja1 the active sheet/selection (previously Hoja3 was active)
# ---------------------------------------------------------------
$ws1.Activate() | Out-Null
$ws1.Range("A1:D4").Select() | Out-Null
